$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.842.87'
$ws.Range("E2").Value = '  +6.18%  '
$ws.Range("D3").Value = '2.760.29'
$ws.Range("E3").Value = '  +4.83%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '118.12'
$ws.Range("E5").Value = '  +7.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '332.43'
$ws.Range("E6").Value = '  +3.18%  '
$ws.Range("E7").Value = '  +3.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.16%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.579'
$ws.Range("E9").Value = '  +7.46%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.09'
$ws.Range("E10").Value = '  +6.81%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.10'
$ws.Range("E11").Value = '  +1.46%  '
$ws.Range("E12").Value = '  +2.76%  '
$ws.Range("E13").Value = '  +3.18%  '
$ws.Range("E14").Value = '  +6.11%  '
$ws.Range("D15").Value = '3.192.53'
$ws.Range("E15").Value = '  +4.62%  '
$ws.Range("D16").Value = '2.752.53'
$ws.Range("E16").Value = '  +4.22%  '
$ws.Range("E17").Value = '  +3.39%  '
$ws.Range("D18").Value = '51.672.61'
$ws.Range("E18").Value = '  +5.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.68'
$ws.Range("E19").Value = '  +6.98%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.87'
$ws.Range("E20").Value = '  +3.29%  '
$ws.Range("B21").Value = 'ImmutableX'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.98'
$ws.Range("E21").Value = '  +3.44%  '
$ws.Range("D22").Value = '0.0₃0963'
$ws.Range("E22").Value = '  +2.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '278.26'
$ws.Range("E23").Value = '  +3.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.81'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("E25").Value = '  +4.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.86'
$ws.Range("E26").Value = '  +2.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.14'
$ws.Range("E27").Value = '  +0.53%  '
$ws.Range("E28").Value = '  +0.11%  '
$ws.Range("E29").Value = '  +2.21%  '
$ws.Range("E30").Value = '  +0.20%  '
$ws.Range("E31").Value = '  +2.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '35.45'
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '50.48'
$ws.Range("E33").Value = '  +2.51%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.60'
$ws.Range("E34").Value = '  +3.53%  '
$ws.Range("E35").Value = '  +4.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '19.18'
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("E37").Value = '  -0.22%  '
$ws.Range("E38").Value = '  +4.52%  '
$ws.Range("E39").Value = '  +2.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.23'
$ws.Range("E40").Value = '  +3.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '130.44'
$ws.Range("E41").Value = '  +4.54%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '23.38'
$ws.Range("E42").Value = '  +4.34%  '
$ws.Range("E43").Value = '  +10.70%  '
$ws.Range("E44").Value = '  +3.11%  '
$ws.Range("E45").Value = '  +4.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.38'
$ws.Range("E46").Value = '  +13.48%  '
$ws.Range("D47").Value = '2.112.42'
$ws.Range("E47").Value = '  +1.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.36'
$ws.Range("E48").Value = '  +4.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.26'
$ws.Range("E49").Value = '  +3.88%  '
$ws.Range("E50").Value = '  +8.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.00'
